$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Cells.Item(2, 2).Value = 463
$ws.Cells.Item(3, 2).Value = 463
$ws.Cells.Item(4, 2).Value = 462
$ws.Cells.Item(5, 2).Value = 461
$ws.Cells.Item(6, 2).Value = 460
$ws.Cells.Item(7, 2).Value = 459
$ws.Cells.Item(8, 2).Value = 458
$ws.Cells.Item(9, 2).Value = 457
$ws.Cells.Item(10, 2).Value = 456
$ws.Cells.Item(11, 2).Value = 456
$ws.Cells.Item(12, 2).Value = 455
$ws.Cells.Item(13, 2).Value = 454
$ws.Cells.Item(14, 2).Value = 453
$ws.Cells.Item(15, 2).Value = 452
$ws.Cells.Item(16, 2).Value = 452
$ws.Cells.Item(17, 2).Value = 451
$ws.Cells.Item(18, 2).Value = 450
$ws.Cells.Item(19, 2).Value = 449
$ws.Cells.Item(20, 2).Value = 448
$ws.Cells.Item(21, 2).Value = 447
$ws.Cells.Item(22, 2).Value = 447
$ws.Cells.Item(23, 2).Value = 446
$ws.Cells.Item(24, 2).Value = 445
$ws.Cells.Item(25, 2).Value = 445
$ws.Cells.Item(26, 2).Value = 444
$ws.Cells.Item(27, 2).Value = 443
$ws.Cells.Item(28, 2).Value = 442
$ws.Cells.Item(29, 2).Value = 440
$ws.Cells.Item(30, 2).Value = 440
$ws.Cells.Item(31, 2).Value = 439
$ws.Cells.Item(32, 2).Value = 438
$ws.Cells.Item(33, 2).Value = 437
$ws.Cells.Item(34, 2).Value = 437
$ws.Cells.Item(35, 2).Value = 435
$ws.Cells.Item(36, 2).Value = 435
$ws.Cells.Item(37, 2).Value = 434
$ws.Cells.Item(38, 2).Value = 433
$ws.Cells.Item(39, 2).Value = 432
$ws.Cells.Item(40, 2).Value = 431
$ws.Cells.Item(41, 2).Value = 430
$ws.Cells.Item(42, 2).Value = 427
$ws.Cells.Item(43, 2).Value = 426
$ws.Cells.Item(44, 2).Value = 426
$ws.Cells.Item(45, 2).Value = 423
$ws.Cells.Item(46, 2).Value = 422
$ws.Cells.Item(47, 2).Value = 421
$ws.Cells.Item(48, 2).Value = 421
$ws.Cells.Item(49, 2).Value = 418
$ws.Cells.Item(50, 2).Value = 417
$ws.Cells.Item(51, 2).Value = 414
$ws.Cells.Item(52, 2).Value = 413
$ws.Cells.Item(53, 2).Value = 411
$ws.Cells.Item(54, 2).Value = 410
$ws.Cells.Item(55, 2).Value = 409
$ws.Cells.Item(56, 2).Value = 408
$ws.Cells.Item(57, 2).Value = 407
$ws.Cells.Item(58, 2).Value = 407
$ws.Cells.Item(59, 2).Value = 405
$ws.Cells.Item(60, 2).Value = 404
$ws.Cells.Item(61, 2).Value = 402
$ws.Cells.Item(62, 2).Value = 401
$ws.Cells.Item(63, 2).Value = 400
$ws.Cells.Item(64, 2).Value = 400
$ws.Cells.Item(65, 2).Value = 399
$ws.Cells.Item(66, 2).Value = 398
$ws.Cells.Item(67, 2).Value = 395
$ws.Cells.Item(68, 2).Value = 387
$ws.Cells.Item(69, 2).Value = 386
$ws.Cells.Item(70, 2).Value = 386
$ws.Cells.Item(71, 2).Value = 385
$ws.Cells.Item(72, 2).Value = 384
$ws.Cells.Item(73, 2).Value = 383
$ws.Cells.Item(74, 2).Value = 383
$ws.Cells.Item(75, 2).Value = 382
$ws.Cells.Item(76, 2).Value = 381
$ws.Cells.Item(77, 2).Value = 381
